$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new "property_category" column before the existing "date" column,
# shifting date / legislator_name / legislator_id one column to the right.
$ws.Columns("H:H").Insert()
$ws.Range("H1").Value = "property_category"
$ws.Range("H2").Value = "stock"
$ws.Range("H3").Value = "stock"
$ws.Range("H4").Value = "stock"

# Fix a typo in the company name (stray space before 司).
$ws.Range("B4").Value = "國票金融控股股份有限公司"
